# Feature/dd20 format change
# Update the two "E_GGG.." shared strings used by the DD20Mapping sheet and
# move the active selection to M3:M4 (matching the refreshed mapping sample).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the example DD20/ETS name pair shown in row 2.
$ws.Range("A2").Value = "E_EEE-FFF_1"
$ws.Range("B2").Value = "E_EEEV-FFF_1"

# Move the current selection/active cell to M3:M4.
$ws.Range("M3:M4").Select()
